$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $escaped = $val -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial(-4163) | Out-Null
}

Set-TextValue 'D2' '67.534.31'
Set-TextValue 'E2' '  +1.32%  '
Set-TextValue 'D3' '3.535.53'
Set-TextValue 'E3' '  +0.99%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '597.61'
Set-TextValue 'E5' '  +1.26%  '
Set-TextValue 'D6' '173.87'
Set-TextValue 'E6' '  +2.68%  '
Set-TextValue 'E7' '  -0.02%  '
Set-TextValue 'D8' '0.594'
Set-TextValue 'E8' '  +2.65%  '
Set-TextValue 'D9' '0.135'
Set-TextValue 'E9' '  +8.74%  '
Set-TextValue 'D10' '7.32'
Set-TextValue 'E10' '  +0.86%  '
Set-TextValue 'E11' '  +0.18%  '
Set-TextValue 'D12' '4.141.45'
Set-TextValue 'E12' '  +0.88%  '
Set-TextValue 'E13' '  -0.16%  '
Set-TextValue 'D14' '28.91'
Set-TextValue 'E14' '  +3.33%  '
Set-TextValue 'D15' '0.0000183'
Set-TextValue 'E15' '  +3.07%  '
Set-TextValue 'D16' '67.417.82'
Set-TextValue 'E16' '  +1.26%  '
Set-TextValue 'D17' '3.535.39'
Set-TextValue 'E17' '  +0.74%  '
Set-TextValue 'E18' '  +1.36%  '
Set-TextValue 'D19' '14.22'
Set-TextValue 'E19' '  +1.81%  '
Set-TextValue 'D20' '397.94'
Set-TextValue 'E20' '  +2.45%  '
Set-TextValue 'D21' '8.01'
Set-TextValue 'E21' '  +0.60%  '
Set-TextValue 'D22' '73.61'
Set-TextValue 'E22' '  +0.83%  '
Set-TextValue 'E23' '  +2.95%  '
Set-TextValue 'D24' '1.00'
Set-TextValue 'E24' '  -0.18%  '
Set-TextValue 'E25' '  +0.74%  '
Set-TextValue 'D26' '10.30'
Set-TextValue 'E26' '  +1.95%  '
Set-TextValue 'E27' '  +0.55%  '
Set-TextValue 'D28' '0.998'
Set-TextValue 'E28' '  -0.16%  '
Set-TextValue 'E29' '  -0.49%  '
Set-TextValue 'D30' '1.48'
Set-TextValue 'E30' '  +1.03%  '
Set-TextValue 'E31' '  +1.60%  '
Set-TextValue 'D32' '24.15'
Set-TextValue 'E32' '  +3.00%  '
Set-TextValue 'D33' '7.41'
Set-TextValue 'E33' '  +0.17%  '
Set-TextValue 'E34' '  +4.60%  '
Set-TextValue 'D35' '163.92'
Set-TextValue 'E35' '  +1.91%  '
Set-TextValue 'E36' '  -0.34%  '
Set-TextValue 'E37' '  -0.34%  '
Set-TextValue 'D38' '6.98'
Set-TextValue 'E38' '  +4.43%  '
Set-TextValue 'D39' '4.75'
Set-TextValue 'E39' '  +2.39%  '
Set-TextValue 'D40' '0.0751'
Set-TextValue 'E40' '  +0.64%  '
Set-TextValue 'E41' '  +0.90%  '
Set-TextValue 'B42' 'InjectiveProtocol'
Set-TextValue 'C42' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D42' '27.37'
Set-TextValue 'E42' '  +2.30%  '
Set-TextValue 'B43' 'dogwifhat'
Set-TextValue 'C43' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D43' '2.65'
Set-TextValue 'E43' '  +4.86%  '
Set-TextValue 'D44' '2.814.74'
Set-TextValue 'E44' '  +0.72%  '
Set-TextValue 'D45' '43.03'
Set-TextValue 'E45' '  -0.84%  '
Set-TextValue 'E46' '  -0.95%  '
Set-TextValue 'D47' '342.92'
Set-TextValue 'E47' '  -3.28%  '
Set-TextValue 'E48' '  +1.13%  '
Set-TextValue 'D49' '33.99'
Set-TextValue 'E49' '  +2.59%  '
Set-TextValue 'D50' '6.54'
Set-TextValue 'E50' '  +0.99%  '
Set-TextValue 'D51' '0.856'
Set-TextValue 'E51' '  +1.04%  '

$excel.CutCopyMode = $false
Write-Host "Applied all cell updates"
